$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 2 (account 005135105 / Brenner / 611705.04),
# shifting all rows below it up by one.
$ws.Rows.Item(2).Delete()
